# Add a new "2022-Q3" sheet (duplicate of the "2022-Q2" sheet, trimmed to
# its own 2 funds) right before the existing "2022-Q2" sheet, and insert a
# matching summary row at the top of the "总计" sheet's data.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Create the new "2022-Q3" worksheet by duplicating "2022-Q2" (keeps
#    headers/number formats/borders identical) and inserting it directly
#    before "2022-Q2".
# ---------------------------------------------------------------------
$q2 = $wb.Worksheets.Item("2022-Q2")
$q2idx = $q2.Index
$q2.Copy($q2)
# The freshly inserted copy lands at $q2idx (taking the "(2)" disambiguated
# name); the original "2022-Q2" sheet is pushed one slot later.
$q3 = $wb.Worksheets.Item($q2idx)
$q3.Name = "2022-Q3"

# Trim the copied sheet from 5 data rows down to the 2 that belong to
# 2022-Q3 (rows 4 and 5 are no longer needed).
$q3.Rows.Item(4).Delete() | Out-Null
$q3.Rows.Item(4).Delete() | Out-Null

# Fill in the 2022-Q3 fund data.
$q3.Range("A2").Value = 0
$q3.Range("B2").Value = "002295"
$q3.Range("C2").Value = "广发稳安灵活配置混合A"
$q3.Range("D2").Value = "'1.58"
$q3.Range("E2").Value = "'69.63"
$q3.Range("F2").Value = "'4.49"
$q3.Range("G2").Value = "'0.0709"
$q3.Range("H2").Value = 3

$q3.Range("A3").Value = 1
$q3.Range("B3").Value = "008604"
$q3.Range("C3").Value = "广发稳安灵活配置混合C"
$q3.Range("D3").Value = "'0.02"
$q3.Range("E3").Value = "'69.63"
$q3.Range("F3").Value = "'4.49"
$q3.Range("G3").Value = "'0.0009"
$q3.Range("H3").Value = 3

# ---------------------------------------------------------------------
# 2. Insert the 2022-Q3 summary row at the top of the data in "总计".
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")
$total.Rows.Item(2).Insert() | Out-Null

$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q3"
$total.Range("C2").Value = 2
$total.Range("D2").Value = 0.07

# Renumber the 0-based index column (A) for the rows that shifted down.
$total.Range("A3").Value = 1
$total.Range("A4").Value = 2
$total.Range("A5").Value = 3
$total.Range("A6").Value = 4
$total.Range("A7").Value = 5
$total.Range("A8").Value = 6
$total.Range("A9").Value = 7
